$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (also reflected in the workbook's sheet list)
$ws.Name = "batch_size_1.0"

# Update row 2 values
$ws.Range("B2").Value = 0.5070674419403076
$ws.Range("E2").Value = 0.1021833333333333
$ws.Range("F2").Value = 0.101

# Update row 3 values
$ws.Range("B3").Value = 0.5948929786682129
$ws.Range("C3").Value = 20.28106538184927
$ws.Range("E3").Value = 0.9706166666666667
$ws.Range("F3").Value = 0.9752999999999999
